$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the "last updated" date (C1) ---
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = 45392

# --- "MCF" sheet: raise capacity factors to 1 for the affected plant types ---
$mcf = $wb.Worksheets.Item("MCF")

# Hard-coded capacity-factor cells (rows whose B value is a literal number,
# not a formula) that moved from 0.85/0.95 up to 1.
$rows = @(2,3,4,6,10,11,12,13,14,16,17,18)
foreach ($r in $rows) {
    $mcf.Range("B$r").Value = 1
}

# Rows 19,20,21,22,24,25 hold formulas (=B2, =B4, =B10, =B14, =B4, =B4) that
# reference the cells above - leave the formulas intact; their cached
# results recalculate to 1 automatically once the precedents change.

# Move the active selection on the MCF sheet to B17 (matches saved view state)
$mcf.Range("B17").Select()
